# Commit: "added a newmain file"
# Append " (Changed main)" -- split across three new runs -- right after the
# existing text run(s) of the document's first paragraph
# ("This is a Microsoft word document."), leaving the original run content
# and the paragraph's own attributes (paraId/textId/rsids) untouched.

$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1).Range

function Escape-Xml([string]$s) {
    $s = $s -replace '&', '&amp;'
    $s = $s -replace '<', '&lt;'
    $s = $s -replace '>', '&gt;'
    return $s
}

# Grab the paragraph's existing text (no trailing pilcrow) and the real
# <w:p ...> attributes (w14:paraId / w14:textId / rsids) off the live
# paragraph so the rebuilt <w:p> that InsertXML writes is identical to what
# is already in the document, save for the new runs appended at the end.
# Paragraph.Range.Text carries a trailing paragraph-mark char (CR/LF) that
# is not part of the visible run text, so trim it off.
$origText = $p1.Text.TrimEnd([char]13, [char]10, [char]7)
$owx = $p1.WordOpenXML

$pAttrs = ""
if ($owx -match '<w:p ([^>]*)>') {
    $pAttrs = $matches[1]
}
$pOpenTag = "<w:p>"
if ($pAttrs -ne "") {
    $pOpenTag = "<w:p " + $pAttrs + ">"
}

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' `
     + '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' `
     + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' `
     + '<pkg:xmlData>' `
     + '<w:document ' + $ns + '>' `
     + '<w:body>' `
     + $pOpenTag `
     + '<w:r><w:t>' + (Escape-Xml $origText) + '</w:t></w:r>' `
     + '<w:r><w:t xml:space="preserve"> (</w:t></w:r>' `
     + '<w:r><w:t>Changed main</w:t></w:r>' `
     + '<w:r><w:t>)</w:t></w:r>' `
     + '</w:p>' `
     + '<w:sectPr/>' `
     + '</w:body>' `
     + '</w:document>' `
     + '</pkg:xmlData></pkg:part></pkg:package>'

$p1.InsertXML($xml) | Out-Null
